# Lillekat budget workbook update (sheet "2023" / sheet2.xml):
#
# The small "Lille Kat 2022-F" side-table in columns K:L gets a new income
# line, "Stotte fra DDSA" = 10000 kr., inserted right below "Start
# kapital" (old row 5). Every row of that table from the old row 6 down
# moves one row lower. Columns N:O ("Lille Kat 2022-S") are a separate,
# independent table and stay exactly where they are.
# Also: the "Lille Kat 2022.10" expense line (now row 14) gets an actual
# value of 7182 instead of being blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Udgifter" header (old K9:L9) is merged; unmerge before shuffling
# rows around so the per-row Copy() below isn't fighting a merged range.
$ws.Range("K9:L9").UnMerge()

# Push the K:L block (old rows 6-21) down by one row. Go bottom-up so we
# never overwrite a row before we've copied it. Copy() brings the cell
# style with it, so formatting/borders land exactly where they were.
for ($r = 21; $r -ge 6; $r--) {
    $ws.Range("K" + $r + ":L" + $r).Copy($ws.Range("K" + ($r + 1) + ":L" + ($r + 1)))
}

# Copy() of a blank source cell doesn't blank the destination in this
# engine, so rows that end up with no K:L content in the new layout (the
# old blank-separator rows, now shifted down by one) still carry stale
# values/formatting from before the shift. Clear them explicitly.
$ws.Range("K9:L9").ClearContents()
$ws.Range("K9:L9").ClearFormats()
$ws.Range("K18:L18").ClearContents()
$ws.Range("K18:L18").ClearFormats()
$ws.Range("K21:L21").ClearContents()
$ws.Range("K21:L21").ClearFormats()

# Re-merge the "Udgifter" header at its new home, K10:L10.
$ws.Range("K10:L10").Merge()

# New row 6: "Stotte fra DDSA" / 10000, styled like the "Start kapital"
# row right above it.
$ws.Range("K5:L5").Copy($ws.Range("K6:L6"))
$ws.Range("K6").Value = "Støtte fra DDSA"
$ws.Range("L6").Value = 10000

# Fix up the formulas that the row shift leaves pointing at the wrong
# ranges (Copy() in this engine pastes values, not live formulas, so we
# set these explicitly to match the new layout).
$ws.Range("L8").Formula = "=SUM(L5:L7)"
$ws.Range("L17").Formula = "=SUM(L11:L16)"
$ws.Range("L19").Formula = "=L8-L17"

# "Lille Kat 2022.10" (row 14) now has an actual spend value.
$ws.Range("L14").Value = 7182

# Row height / thick border hints so the visual "double rule" under the
# two new subtotal rows (8 and 17) and the Balance row (19) match.
$ws.Range("A8:O8").RowHeight = 18
$ws.Range("A9:O9").RowHeight = 17
$ws.Range("A10:O10").RowHeight = 17
$ws.Range("A11:O11").RowHeight = 18
$ws.Range("A12:O12").RowHeight = 17
$ws.Range("A17:O17").RowHeight = 18
$ws.Range("A18:O18").RowHeight = 18
$ws.Range("A19:O19").RowHeight = 18
$ws.Range("A20:O20").RowHeight = 17

Write-Host "edit complete"
